$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/15/2024  Through  7/21/2024"

# --- Row 14 (Murder) ---
$ws.Range("N14").Value = -90.909090909090

# --- Row 15 (Rape) ---
$ws.Range("M15").Value = 50

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 117
$ws.Range("J16").Value = 119
$ws.Range("K16").Value = -1.680672268907
$ws.Range("L16").Value = 2.631578947368
$ws.Range("M16").Value = -14.598540145985
$ws.Range("N16").Value = -77.5

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 174
$ws.Range("J17").Value = 163
$ws.Range("K17").Value = 6.748466257668
$ws.Range("L17").Value = 9.433962264150
$ws.Range("M17").Value = 38.095238095238
$ws.Range("N17").Value = 5.454545454545

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 24
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 14
$ws.Range("L18").Value = -22.448979591836
$ws.Range("M18").Value = -55.46875
$ws.Range("N18").Value = -90.008764241893

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -13.461538461538
$ws.Range("I19").Value = 377
$ws.Range("J19").Value = 370
$ws.Range("K19").Value = 1.891891891891
$ws.Range("L19").Value = 7.714285714285
$ws.Range("M19").Value = 62.5
$ws.Range("N19").Value = 10.233918128655

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -10
$ws.Range("I20").Value = 214
$ws.Range("J20").Value = 187
$ws.Range("K20").Value = 14.438502673796
$ws.Range("L20").Value = 46.575342465753
$ws.Range("M20").Value = 5.418719211822
$ws.Range("N20").Value = -89.316025961058

# --- Row 21 (TOTAL) ---
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -17.073170731707
$ws.Range("F21").Value = 146
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 3.546099290780
$ws.Range("I21").Value = 1012
$ws.Range("J21").Value = 952
$ws.Range("K21").Value = 6.302521008403
$ws.Range("L21").Value = 8.934337997847
$ws.Range("M21").Value = 4.870466321243
$ws.Range("N21").Value = -75.899023577042

# --- Row 22 (Transit) : D22 and E22 become the "N/A" text markers ---
# (use a cell-copy so the result is a genuine shared-string "text" cell,
#  not a numeric cell auto-parsed from a numeric-looking string)
$ws.Range("C22").Copy($ws.Range("D22"))
$ws.Range("N22").Copy($ws.Range("E22"))

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -6.25
$ws.Range("F24").Value = 133
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = 29.126213592233
$ws.Range("I24").Value = 799
$ws.Range("J24").Value = 712
$ws.Range("K24").Value = 12.219101123595
$ws.Range("L24").Value = 1.783439490445
$ws.Range("M24").Value = 32.066115702479

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 55
$ws.Range("H25").Value = 30.952380952381
$ws.Range("I25").Value = 294
$ws.Range("J25").Value = 241
$ws.Range("K25").Value = 21.991701244813
$ws.Range("L25").Value = -11.711711711711

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 9
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 12.5
$ws.Range("I26").Value = 341
$ws.Range("J26").Value = 277
$ws.Range("K26").Value = 23.104693140794
$ws.Range("L26").Value = 16.780821917808
$ws.Range("M26").Value = -18.421052631578

# --- Row 27 (UCR Rape*) ---
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 5.263157894736
$ws.Range("L27").Value = 5.263157894736

# --- Row 28 (Other Sex Crimes): C28 becomes number, D28/E28 become "N/A" text markers ---
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 1
$ws.Range("C22").Copy($ws.Range("D28"))
$ws.Range("N22").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -56.818181818181

# --- Row 29 (Shooting Vic.) ---
$ws.Range("N29").Value = -94.117647058823

# --- Row 30 (Shooting Inc.) ---
$ws.Range("N30").Value = -93.75

# --- Row 31 (Hate Crimes): C31 and F31 become numbers ---
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("C31").Value = 1
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 5
$ws.Range("L31").Value = 25
